$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell text values (content changes), in the order the new shared
# strings appear in the diff: C10, C9, A7, A14
$ws.Range("C10").Value = "pep[c] ---> pep[e]"
$ws.Range("C9").Value = "enz[c] ---> enz[e]"
$ws.Range("A7").Value = "bmt2r"
$ws.Range("A14").Value = "bmex"

# Update the selection to A14 as the active cell
$ws.Range("A14").Select()
